$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing year column (R, rows 2-5) into
# the new column S so the new cells inherit the same borders/number
# formats/fonts as the rest of the table. Row 1 (the merged title band)
# does not get a new cell in column S, so it is intentionally excluded.
$ws.Range("R2:R5").Copy()
$ws.Range("S2:S5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New year column: 2022, with its population figure and the (unchanged)
# percentage-of-population ratio repeated from the previous year.
$ws.Range("S3").Value = 2022
$ws.Range("S4").Value = 211650
$ws.Range("S5").Value = 2.9794303052841493

# Move the active selection to the new last cell of row 2, matching where
# a user would land after extending the table by one more column.
[void]$ws.Range("S2").Select()
